$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 41: FraxShare -> VeChain (rows 41/42 swapped order)
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D41") "0.02354"
$ws.Range("E41").Value = "  +1.90%  "

# Row 42: VeChain -> FraxShare
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D42") "8.821"
$ws.Range("E42").Value = "  +5.27%  "

# Remaining rows: update Price (D) and Volume(1h) (E) only
Set-TextValue $ws.Range("D2") "27.393.86"
$ws.Range("E2").Value = "  +3.49%  "
Set-TextValue $ws.Range("D3") "1.797.52"
$ws.Range("E3").Value = "  +4.63%  "
Set-TextValue $ws.Range("D4") "1.005"
$ws.Range("E4").Value = "  +0.11%  "
Set-TextValue $ws.Range("D5") "335.94"
$ws.Range("E5").Value = "  +0.83%  "
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  -0.02%  "
Set-TextValue $ws.Range("D7") "0.3786"
$ws.Range("E7").Value = "  +2.25%  "
Set-TextValue $ws.Range("D8") "48.88"
$ws.Range("E8").Value = "  +1.40%  "
Set-TextValue $ws.Range("D9") "0.3444"
$ws.Range("E9").Value = "  +2.82%  "
Set-TextValue $ws.Range("D10") "1.208"
$ws.Range("E10").Value = "  +2.14%  "
Set-TextValue $ws.Range("D11") "0.07496"
$ws.Range("E11").Value = "  +1.72%  "
Set-TextValue $ws.Range("D12") "1.002"
$ws.Range("E12").Value = "  +0.02%  "
Set-TextValue $ws.Range("D13") "22.08"
$ws.Range("E13").Value = "  +10.22%  "
Set-TextValue $ws.Range("D14") "6.507"
$ws.Range("E14").Value = "  +2.33%  "
Set-TextValue $ws.Range("D15") "1.789.60"
$ws.Range("E15").Value = "  +4.07%  "
Set-TextValue $ws.Range("D16") "7.061"
$ws.Range("E16").Value = "  +0.59%  "
Set-TextValue $ws.Range("D17") "0.00001100"
$ws.Range("E17").Value = "  +3.22%  "
Set-TextValue $ws.Range("D18") "0.06646"
$ws.Range("E18").Value = "  +0.36%  "
Set-TextValue $ws.Range("D19") "84.72"
$ws.Range("E19").Value = "  +3.58%  "
Set-TextValue $ws.Range("D20") "1.002"
$ws.Range("E20").Value = "  +0.15%  "
Set-TextValue $ws.Range("D21") "17.46"
$ws.Range("E21").Value = "  +5.94%  "
Set-TextValue $ws.Range("D22") "6.512"
$ws.Range("E22").Value = "  +6.59%  "
Set-TextValue $ws.Range("D23") "27.407.37"
$ws.Range("E23").Value = "  +3.59%  "
Set-TextValue $ws.Range("D24") "12.56"
$ws.Range("E24").Value = "  -1.35%  "
Set-TextValue $ws.Range("D25") "2.462"
$ws.Range("E25").Value = "  +1.39%  "
Set-TextValue $ws.Range("D26") "1.548"
$ws.Range("E26").Value = "  +12.04%  "
Set-TextValue $ws.Range("D27") "2.584"
$ws.Range("E27").Value = "  +8.64%  "
Set-TextValue $ws.Range("D28") "21.59"
$ws.Range("E28").Value = "  +11.80%  "
Set-TextValue $ws.Range("D29") "150.71"
$ws.Range("E29").Value = "  -0.39%  "
Set-TextValue $ws.Range("D30") "1.989.67"
$ws.Range("E30").Value = "  +4.07%  "
Set-TextValue $ws.Range("D31") "133.81"
$ws.Range("E31").Value = "  +2.26%  "
Set-TextValue $ws.Range("D32") "4.075"
$ws.Range("E32").Value = "  -1.00%  "
Set-TextValue $ws.Range("D33") "6.178"
$ws.Range("E33").Value = "  +4.89%  "
Set-TextValue $ws.Range("D34") "0.08655"
$ws.Range("E34").Value = "  +0.59%  "
Set-TextValue $ws.Range("D35") "13.30"
$ws.Range("E35").Value = "  +5.48%  "
Set-TextValue $ws.Range("D36") "1.686"
$ws.Range("E36").Value = "  -0.93%  "
Set-TextValue $ws.Range("D37") "5.473"
$ws.Range("E37").Value = "  +2.62%  "
Set-TextValue $ws.Range("D38") "0.6924"
$ws.Range("E38").Value = "  +12.50%  "
Set-TextValue $ws.Range("D39") "0.2212"
$ws.Range("E39").Value = "  +2.87%  "
Set-TextValue $ws.Range("D40") "0.06372"
$ws.Range("E40").Value = "  +3.15%  "
Set-TextValue $ws.Range("D43") "1.275"
$ws.Range("E43").Value = "  +4.80%  "
Set-TextValue $ws.Range("D44") "14.57"
$ws.Range("E44").Value = "  +2.96%  "
Set-TextValue $ws.Range("D45") "0.6495"
$ws.Range("E45").Value = "  +9.15%  "
Set-TextValue $ws.Range("D46") "1.003"
$ws.Range("E46").Value = "  +0.22%  "
Set-TextValue $ws.Range("D47") "3.847"
$ws.Range("E47").Value = "  -1.40%  "
Set-TextValue $ws.Range("D48") "2.135"
$ws.Range("E48").Value = "  +5.18%  "
Set-TextValue $ws.Range("D49") "130.27"
$ws.Range("E49").Value = "  +1.68%  "
Set-TextValue $ws.Range("D50") "0.07207"
$ws.Range("E50").Value = "  +0.71%  "
Set-TextValue $ws.Range("D51") "79.58"
$ws.Range("E51").Value = "  +3.84%  "
